# Weekly update: insert a new week's price row for
# "Hortaliza, Feria Lagunitas de Puerto Montt - Poroto verde" at row 47,
# pushing the existing rows 47-96 down to 48-97, and populate the new
# row 47 with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 47 (shifts 47..96 -> 48..97)
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with this week's record.
$ws.Cells.Item(47, 1).Value = 4
$ws.Cells.Item(47, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(47, 3).Value = "Los Lagos"
$ws.Cells.Item(47, 4).Value = 44803
$ws.Cells.Item(47, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(47, 5).Value = 10
$ws.Cells.Item(47, 6).Value = 100112031
$ws.Cells.Item(47, 7).Value = "Poroto verde"
$ws.Cells.Item(47, 8).Value = "Magnum"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 40
$ws.Cells.Item(47, 11).Value = 36000
$ws.Cells.Item(47, 12).Value = 36000
$ws.Cells.Item(47, 13).Value = 36000
$ws.Cells.Item(47, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(47, 15).Value = "Perú"
$ws.Cells.Item(47, 16).Value = 1440
$ws.Cells.Item(47, 17).Value = 25
$ws.Cells.Item(47, 18).Value = "Hortaliza"
